$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 (B5:D5) was blank; fill it in with the same three color names used
# throughout the sheet, reusing the existing formatting (fill/font/border)
# already applied to matching color cells elsewhere on the sheet, so no new
# styles are introduced.
#   B5 -> BLANCO (same look as B4)
#   C5 -> AZUL   (same look as B2)
#   D5 -> NEGRO  (same look as C2)

$ws.Range("B4").Copy($ws.Range("B5"))
$ws.Range("B2").Copy($ws.Range("C5"))
$ws.Range("C2").Copy($ws.Range("D5"))

# Move the active selection to D5.
$ws.Range("D5").Select()
